$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g18.3")

$ws.Range("B22").Value = 22.63496935893463
$ws.Range("C22").Value = 11.48367918335759
$ws.Range("D22").Value = 13.5036084021264

$ws.Range("B23").Value = 23.10012769267524
$ws.Range("C23").Value = 11.8369437163957
$ws.Range("D23").Value = 13.73462385747906

$ws.Range("B24").Value = 23.52642530643275
$ws.Range("C24").Value = 12.16335479852913
$ws.Range("D24").Value = 13.98063481830828

$ws.Range("B25").Value = 23.71046689484384
$ws.Range("C25").Value = 12.45142418607468
$ws.Range("D25").Value = 14.00717120979344

$ws.Range("B26").Value = 23.91614243923391
$ws.Range("C26").Value = 12.62570199784373
$ws.Range("D26").Value = 14.15440654437449
